$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$data = @(
    @(18, 18),
    @(9, 9),
    @(1, 4),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 7),
    @(1, 6),
    @(1, 7),
    @(1, 7),
    @(1, 3),
    @(1, 3),
    @(1, 7),
    @(1, 6),
    @(1, 7),
    @(1, 7),
    @(1, 6),
    @(1, 8),
    @(1, 7),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(1, 7),
    @(1, 6),
    @(1, 5),
    @(1, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
